# "Dataframe ST.xlsx" update — adds the new "25-sep" day column to Sheet1
# and refreshes the Sheet3 lookup table (and its dependent VLOOKUP column)
# with that day's figures.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1) Sheet1: new date column CA ("25-sep") with one value per product row
#    NOTE: write the header first so "25-sep" lands in the shared-string
#    table before the new Sheet3 product label below.
# ---------------------------------------------------------------------
$ws1.Range("CA1").Value = "25-sep"

$newDayValues = @(
    0,                      # row 2  - LAYS CLASICAS 45GRX54
    9.2380510078173028,     # row 3  - LAYS CLASICAS 94GRX25
    11.860409477431915,     # row 4  - LAYS CLASICAS 145GRX18
    15.792991420179193,     # row 5  - LAYS CLASICAS 249GRX14
    0,                      # row 6  - DORITOS QUESO 45GRX70
    0.17288193666060836,    # row 7  - DORITOS QUESO 85GRX26
    12.940041072508262,     # row 8  - DORITOS QUESO 140GRX19
    11.614763627804068,     # row 9  - PEHUAMAR PAPA LISA 520GX9
    14.982430051219202,     # row 10 - PEHUAMAR ACANALADA 520GX9
    8.5085416773233966,     # row 11 - PEHUAMAR MAICITOS 285GX10
    0,                      # row 12 - 3DMEGAQUESO 95GX24
    4.9773522073788445,     # row 13 - CHEETOS 94GRX24
    0,                      # row 14 - QUAKER AVENA INSTANT 18X380G
    0,                      # row 15 - LAYS CEBOLLA CARAMELIZADA 41GX54
    6.7563355521713309,     # row 16 - LAYS CEBOLLA CARAMELIZADA 85GX25
    0,                      # row 17 - LAYS QUESO Y PIMIENTA 41GX54
    0                       # row 18 - LAYS QUESO Y PIMIENTA 85GX25
)

for ($i = 0; $i -lt $newDayValues.Length; $i++) {
    $ws1.Cells.Item($i + 2, 79).Value = $newDayValues[$i]
}

# ---------------------------------------------------------------------
# 2) Sheet3: lookup table rows 20:36 refreshed for the new day, including
#    the relabeled product in row 24. Column C (VLOOKUP against
#    $A$20:$B$36) recalculates automatically once the table changes.
# ---------------------------------------------------------------------
$ws3.Range("A24").Value = "DORITOS QUESO 85GrX26"

$lookupValues = @(
    16.25762017271817,      # row 20 - 3D QUESO 92GX27
    4.9773522073788445,     # row 21 - CHEETOS 94GRX24
    12.940041072508262,     # row 22 - DORITOS QUESO 140GRX19
    1.0266851514853432,     # row 23 - DORITOS QUESO 40GRX58X1 CH
    0.17288193666060836,    # row 24 - DORITOS QUESO 85GrX26
    6.7563355521713309,     # row 25 - LAYS CEBOLLA CARAMELIZADA 85GX25
    11.860409477431915,     # row 26 - LAYS CLASICAS 145GRX18
    15.792991420179193,     # row 27 - LAYS CLASICAS 249GRX14
    4.4536274889009215,     # row 28 - LAYS CLASICAS 40GX68
    9.2380510078173028,     # row 29 - LAYS CLASICAS 94GRX25
    16.54929098162831,      # row 30 - LAYS ONDAS FH 30GX72
    10.28231674603585,      # row 31 - LAYS ONDAS FH 70GX28
    7.7622672353493476,     # row 32 - LAYS QSO Y CEBOLLA 34GX72
    14.982430051219202,     # row 33 - PEHUAMAR ACANALADA 520GX9
    8.5085416773233966,     # row 34 - PEHUAMAR MAICITOS 285GX10
    11.614763627804068,     # row 35 - PEHUAMAR PAPA LISA 520GX9
    23.133381693891206      # row 36 - QUAKER AVENA INSTANT FORTIF 18X280G
)

for ($i = 0; $i -lt $lookupValues.Length; $i++) {
    $ws3.Cells.Item(20 + $i, 2).Value = $lookupValues[$i]
}

# ---------------------------------------------------------------------
# 3) Restore the selections recorded in each sheet's view, ending back on
#    Sheet1 so it remains the active tab.
# ---------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("C7:C8").Select()

$ws1.Activate()
$ws1.Range("CC8").Select()
